# This workbook holds a weekly price table for "Coliflor" (Macroferia Regional
# de Talca). A new week of data (2 rows, dated 44830) is inserted at the top
# of the data block (rows 216-217), which pushes every subsequent row down by
# two positions. The two rows that fall off the bottom of the original block
# (old rows 343 and 344) are appended as new rows 345 and 346 at the end of
# the sheet.
#
# Only these columns actually vary row to row: D (Fecha), I (Calidad),
# J (Volumen), K (Precio minimo), L (Precio maximo), M (Precio promedio
# ponderado), O (Origen), P (Precio $/Kg). All the other columns
# (A,B,C,E,F,G,H,N,Q,R) are constant for every data row in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 216
$lastDataRow = 344
$shiftBy = 2

$varCols = @("D","I","J","K","L","M","O","P")
$fixedCols = @("A","B","C","E","F","G","H","N","Q","R")

# --- 1. Capture the rows that will be pushed past the end of the table ---
# old row (lastDataRow - 1) -> new row (lastDataRow + 1)
# old row (lastDataRow)     -> new row (lastDataRow + 2)
$tailRows = @($lastDataRow - 1, $lastDataRow)
$capturedVar = @{}
foreach ($r in $tailRows) {
    $capturedVar[$r] = @{}
    foreach ($c in $varCols) {
        $capturedVar[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Fixed columns are identical on every data row, so grab one full template row.
$fixedVals = @{}
foreach ($c in $fixedCols) {
    $fixedVals[$c] = $ws.Range("$c$firstDataRow").Value2
}

# --- 2. Write the two brand-new appended rows (old 343 -> 345, old 344 -> 346) ---
$newRow1 = $lastDataRow + 1
$newRow2 = $lastDataRow + 2

foreach ($c in $fixedCols) {
    $ws.Range("$c$newRow1").Value = $fixedVals[$c]
    $ws.Range("$c$newRow2").Value = $fixedVals[$c]
}
foreach ($c in $varCols) {
    $ws.Range("$c$newRow1").Value = $capturedVar[$tailRows[0]][$c]
    $ws.Range("$c$newRow2").Value = $capturedVar[$tailRows[1]][$c]
}
# Date column needs the same date-style number format as the rest of column D.
$ws.Range("D$newRow1").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D$newRow2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- 3. Shift rows (firstDataRow+2) .. lastDataRow down by 2 rows ---
# Walk from the bottom up so we never overwrite a source row before reading it.
for ($r = $lastDataRow; $r -ge ($firstDataRow + $shiftBy); $r--) {
    $src = $r - $shiftBy
    foreach ($c in $varCols) {
        $ws.Range("$c$r").Value = $ws.Range("$c$src").Value2
    }
}

# --- 4. Fill in the brand-new first two rows with this week's data ---
$ws.Range("D$firstDataRow").Value = 44830
$ws.Range("J$firstDataRow").Value = 2000
$ws.Range("K$firstDataRow").Value = 1500
$ws.Range("L$firstDataRow").Value = 1500
$ws.Range("M$firstDataRow").Value = 1500
$ws.Range("P$firstDataRow").Value = 1500
# I216 stays "Primera" and O216 stays "Region del Maule" - no change required.

$secondRow = $firstDataRow + 1
$ws.Range("D$secondRow").Value = 44830
$ws.Range("I$secondRow").Value = "Segunda"
$ws.Range("J$secondRow").Value = 2000
$ws.Range("K$secondRow").Value = 1200
$ws.Range("L$secondRow").Value = 1200
$ws.Range("M$secondRow").Value = 1200
$ws.Range("P$secondRow").Value = 1200
